# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Source data only ever downloaded team statistics; these three new
# columns carry the season record (wins/losses/ties) for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (bold, centered, bordered) from the
# last header cell (AC1) onto the three new header cells so they match the
# rest of row 1 exactly.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record is constant for every player on this roster (they all
# belong to the same team/season): 65 wins, 97 losses, 0 ties.
$wins = 65
$losses = 97
$ties = 0

$lastRow = 57
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}
